$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.056.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.66%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.245.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'396.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'108.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.53%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +7.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.242.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.37%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.33%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'39.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.23%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +9.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +2.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.759.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'8.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.82%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'19.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.358.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +5.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -2.72%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.18%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'56.954.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.86%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +7.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.16%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'293.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +6.81%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'74.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -3.08%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'28.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.03%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.30%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Kaspa"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'0.169"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'RenderToken"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'7.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.96%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.72%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'40.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +9.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0489"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'51.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.04%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -4.58%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.52%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'139.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.62%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.89%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -1.66%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'NEARProtocol"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.37%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Celestia"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'17.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.34%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -4.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +11.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.158.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.71%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -5.78%  "
$ws.Range("E51").Style = "Normal"
